$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.195.06"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "1.960.10"
$ws.Range("E3").Value = "  +2.36%  "

$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'247.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").Value = "'0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "'0.4893"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.06%  "

$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'44.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2972"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.83%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06840"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'19.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "'106.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.79%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07770"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.932.75"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.424"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.7118"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.19%  "

$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "'287.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.22%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "31.202.55"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("D19").Value = "'0.000007773"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.82%  "

$ws.Range("B21").Value = "BitDAO"
$ws.Range("C21").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D21").Value = "'0.4885"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.80%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.591"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'0.9995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.181.37"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "'0.9934"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("B26").Value = "Chainlink"
$ws.Range("C26").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D26").Value = "'6.593"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.75%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.63%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'168.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.11%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'2.205"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.15%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1066"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.441"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.783"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +18.21%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.498"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.09%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.05038"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7689"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.77%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.166"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.28%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02052"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.727"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.714"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.130"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.48%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.14%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8863"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.14%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'109.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'73.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.00%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4458"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.65%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'0.9986"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.94%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'992.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.44%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1271"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.71%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.406"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.70%  "
